# USAID WCS analysis for SPR
# Fix Species names in column K that contained a stray non-breaking
# space (U+00A0) between genus and species; replace with a normal
# ASCII space so the text reads correctly (e.g. "Tylosurus crocodilus").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fixes = @(
    @{ Row = 21; Name = "Tylosurus crocodilus" },
    @{ Row = 22; Name = "Tylosurus crocodilus" },
    @{ Row = 90; Name = "Carangoides ciliarius" },
    @{ Row = 91; Name = "Carangoides ciliarius" },
    @{ Row = 92; Name = "Carangoides ciliarius" },
    @{ Row = 97; Name = "Carangoides malabaricus" },
    @{ Row = 98; Name = "Carangoides malabaricus" },
    @{ Row = 99; Name = "Carangoides malabaricus" },
    @{ Row = 100; Name = "Carangoides malabaricus" },
    @{ Row = 101; Name = "Carangoides malabaricus" },
    @{ Row = 102; Name = "Carangoides malabaricus" },
    @{ Row = 103; Name = "Carangoides malabaricus" },
    @{ Row = 104; Name = "Carangoides malabaricus" },
    @{ Row = 105; Name = "Carangoides malabaricus" },
    @{ Row = 106; Name = "Carangoides malabaricus" },
    @{ Row = 107; Name = "Carangoides malabaricus" },
    @{ Row = 108; Name = "Carangoides malabaricus" },
    @{ Row = 109; Name = "Carangoides malabaricus" },
    @{ Row = 110; Name = "Carangoides malabaricus" },
    @{ Row = 111; Name = "Carangoides malabaricus" },
    @{ Row = 112; Name = "Carangoides malabaricus" },
    @{ Row = 113; Name = "Carangoides malabaricus" },
    @{ Row = 114; Name = "Carangoides malabaricus" },
    @{ Row = 115; Name = "Carangoides malabaricus" },
    @{ Row = 116; Name = "Carangoides malabaricus" },
    @{ Row = 117; Name = "Carangoides malabaricus" },
    @{ Row = 118; Name = "Carangoides malabaricus" },
    @{ Row = 119; Name = "Carangoides malabaricus" },
    @{ Row = 120; Name = "Carangoides malabaricus" },
    @{ Row = 121; Name = "Carangoides malabaricus" },
    @{ Row = 122; Name = "Carangoides malabaricus" },
    @{ Row = 123; Name = "Carangoides malabaricus" },
    @{ Row = 124; Name = "Decapterus kurroides" },
    @{ Row = 125; Name = "Decapterus kurroides" },
    @{ Row = 126; Name = "Decapterus kurroides" },
    @{ Row = 130; Name = "Scomberoides tol" },
    @{ Row = 131; Name = "Scomberoides tol" },
    @{ Row = 132; Name = "Scomberoides tol" },
    @{ Row = 133; Name = "Scomberoides tol" },
    @{ Row = 134; Name = "Scomberoides tol" },
    @{ Row = 135; Name = "Scomberoides tol" },
    @{ Row = 136; Name = "Scomberoides tol" },
    @{ Row = 137; Name = "Scomberoides tol" },
    @{ Row = 138; Name = "Scomberoides tol" },
    @{ Row = 171; Name = "Platax teira" },
    @{ Row = 172; Name = "Platax teira" },
    @{ Row = 173; Name = "Platax teira" },
    @{ Row = 272; Name = "Hyporhamphus affinis" },
    @{ Row = 273; Name = "Hyporhamphus affinis" },
    @{ Row = 274; Name = "Hyporhamphus affinis" },
    @{ Row = 275; Name = "Hyporhamphus affinis" },
    @{ Row = 276; Name = "Hyporhamphus affinis" },
    @{ Row = 277; Name = "Hyporhamphus affinis" },
    @{ Row = 278; Name = "Hyporhamphus affinis" },
    @{ Row = 279; Name = "Hyporhamphus affinis" },
    @{ Row = 280; Name = "Hyporhamphus affinis" },
    @{ Row = 281; Name = "Hyporhamphus affinis" },
    @{ Row = 282; Name = "Hyporhamphus affinis" },
    @{ Row = 283; Name = "Hyporhamphus affinis" },
    @{ Row = 284; Name = "Hyporhamphus affinis" },
    @{ Row = 285; Name = "Hyporhamphus affinis" },
    @{ Row = 286; Name = "Hyporhamphus affinis" },
    @{ Row = 287; Name = "Hyporhamphus affinis" },
    @{ Row = 288; Name = "Hyporhamphus affinis" },
    @{ Row = 289; Name = "Hyporhamphus affinis" },
    @{ Row = 290; Name = "Hyporhamphus affinis" },
    @{ Row = 301; Name = "Cheilinus trilobatus" },
    @{ Row = 302; Name = "Cheilinus trilobatus" },
    @{ Row = 303; Name = "Cheilinus trilobatus" },
    @{ Row = 304; Name = "Cheilinus trilobatus" },
    @{ Row = 550; Name = "Monodactylus argenteus" },
    @{ Row = 551; Name = "Monodactylus argenteus" },
    @{ Row = 552; Name = "Mulloidichthys vanicolensis" },
    @{ Row = 553; Name = "Mulloidichthys vanicolensis" },
    @{ Row = 554; Name = "Mulloidichthys vanicolensis" },
    @{ Row = 625; Name = "Upeneus sulphureus" },
    @{ Row = 626; Name = "Upeneus vittatus" },
    @{ Row = 627; Name = "Upeneus vittatus" },
    @{ Row = 628; Name = "Upeneus vittatus" },
    @{ Row = 629; Name = "Upeneus vittatus" },
    @{ Row = 630; Name = "Upeneus vittatus" },
    @{ Row = 631; Name = "Upeneus vittatus" },
    @{ Row = 632; Name = "Upeneus vittatus" },
    @{ Row = 633; Name = "Upeneus vittatus" },
    @{ Row = 634; Name = "Upeneus vittatus" },
    @{ Row = 635; Name = "Upeneus vittatus" },
    @{ Row = 636; Name = "Upeneus vittatus" },
    @{ Row = 637; Name = "Upeneus vittatus" },
    @{ Row = 755; Name = "Cephalopholis spiloparaea" },
    @{ Row = 759; Name = "Epinephelus spilotoceps" },
    @{ Row = 760; Name = "Epinephelus tauvina" },
    @{ Row = 761; Name = "Epinephelus tauvina" },
    @{ Row = 906; Name = "Sphyraena obtusata" },
    @{ Row = 907; Name = "Sphyraena obtusata" },
    @{ Row = 908; Name = "Sphyraena obtusata" },
    @{ Row = 909; Name = "Sphyraena obtusata" },
    @{ Row = 910; Name = "Sphyraena obtusata" },
    @{ Row = 911; Name = "Sphyraena obtusata" },
    @{ Row = 912; Name = "Sphyraena obtusata" },
    @{ Row = 913; Name = "Sphyraena obtusata" },
    @{ Row = 914; Name = "Sphyraena obtusata" },
    @{ Row = 915; Name = "Sphyraena obtusata" },
    @{ Row = 916; Name = "Sphyraena obtusata" },
    @{ Row = 917; Name = "Sphyraena obtusata" },
    @{ Row = 918; Name = "Sphyraena obtusata" },
    @{ Row = 919; Name = "Sphyraena obtusata" },
    @{ Row = 920; Name = "Sphyraena obtusata" },
    @{ Row = 921; Name = "Sphyraena obtusata" },
    @{ Row = 922; Name = "Sphyraena obtusata" },
    @{ Row = 923; Name = "Sphyraena obtusata" },
    @{ Row = 924; Name = "Sphyraena obtusata" },
    @{ Row = 925; Name = "Sphyraena obtusata" },
    @{ Row = 926; Name = "Sphyraena obtusata" },
    @{ Row = 927; Name = "Sphyraena obtusata" },
    @{ Row = 928; Name = "Sphyraena obtusata" },
    @{ Row = 929; Name = "Sphyraena obtusata" },
    @{ Row = 930; Name = "Sphyraena obtusata" },
    @{ Row = 931; Name = "Sphyraena obtusata" },
    @{ Row = 932; Name = "Sphyraena obtusata" },
    @{ Row = 933; Name = "Sphyraena obtusata" },
    @{ Row = 934; Name = "Sphyraena obtusata" },
    @{ Row = 935; Name = "Sphyraena obtusata" },
    @{ Row = 936; Name = "Sphyraena obtusata" },
    @{ Row = 937; Name = "Sphyraena obtusata" },
    @{ Row = 938; Name = "Sphyraena obtusata" },
    @{ Row = 939; Name = "Sphyraena obtusata" },
    @{ Row = 940; Name = "Sphyraena obtusata" },
    @{ Row = 941; Name = "Sphyraena obtusata" },
    @{ Row = 942; Name = "Sphyraena obtusata" },
    @{ Row = 943; Name = "Sphyraena obtusata" },
    @{ Row = 944; Name = "Sphyraena obtusata" },
    @{ Row = 945; Name = "Sphyraena obtusata" },
    @{ Row = 946; Name = "Sphyraena obtusata" },
    @{ Row = 947; Name = "Sphyraena obtusata" },
    @{ Row = 948; Name = "Sphyraena obtusata" },
    @{ Row = 949; Name = "Sphyraena obtusata" },
    @{ Row = 950; Name = "Sphyraena obtusata" },
    @{ Row = 951; Name = "Sphyraena obtusata" },
    @{ Row = 952; Name = "Sphyraena obtusata" },
    @{ Row = 953; Name = "Sphyraena obtusata" },
    @{ Row = 954; Name = "Sphyraena obtusata" },
    @{ Row = 955; Name = "Saurida gracilis" }
)

foreach ($fix in $fixes) {
    $ws.Cells.Item($fix.Row, 11).Value = $fix.Name
}
